$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 22 de Abril de 2020 a las 02:22"

# Refresh country case-count data and re-sort affected rows (Sudan, Libia, Belice/Malaui, Yemen moves)
# Row 4: Estados Unidos -> Estados Unidos
$ws.Cells.Item(4,2).Value = 819620
$ws.Cells.Item(4,3).Value = 26861
$ws.Cells.Item(4,5).Value = 691446
$ws.Cells.Item(4,7).Value = 2800
$ws.Cells.Item(4,8).Value = 45314
# Row 14: Brasil -> Brasil
$ws.Cells.Item(14,4).Value = 24325
$ws.Cells.Item(14,5).Value = 16013
# Row 38: Noruega -> Noruega
$ws.Cells.Item(38,2).Value = 7241
$ws.Cells.Item(38,3).Value = 85
$ws.Cells.Item(38,5).Value = 7027
# Row 57: Argentina -> Argentina
$ws.Cells.Item(57,2).Value = 3144
$ws.Cells.Item(57,3).Value = 113
$ws.Cells.Item(57,5).Value = 2153
$ws.Cells.Item(57,7).Value = 9
$ws.Cells.Item(57,8).Value = 151
# Row 134: Brunei -> Sudan
$ws.Cells.Item(134,1).Value = "Sudan"
$ws.Cells.Item(134,2).Value = 140
$ws.Cells.Item(134,3).Value = 33
$ws.Cells.Item(134,4).Value = 8
$ws.Cells.Item(134,5).Value = 119
$ws.Cells.Item(134,6).Value = 0
$ws.Cells.Item(134,7).Value = 1
$ws.Cells.Item(134,8).Value = 13
# Row 135: Gibraltar -> Brunei
$ws.Cells.Item(135,1).Value = "Brunei"
$ws.Cells.Item(135,2).Value = 138
$ws.Cells.Item(135,4).Value = 116
$ws.Cells.Item(135,5).Value = 21
$ws.Cells.Item(135,6).Value = 2
$ws.Cells.Item(135,8).Value = 1
# Row 136: Camboya -> Gibraltar
$ws.Cells.Item(136,1).Value = "Gibraltar"
$ws.Cells.Item(136,2).Value = 132
$ws.Cells.Item(136,4).Value = 120
# Row 137: Birmania -> Camboya
$ws.Cells.Item(137,1).Value = "Camboya"
$ws.Cells.Item(137,2).Value = 122
$ws.Cells.Item(137,3).Value = 0
$ws.Cells.Item(137,4).Value = 110
$ws.Cells.Item(137,5).Value = 12
$ws.Cells.Item(137,6).Value = 1
$ws.Cells.Item(137,8).Value = 0
# Row 138: Madagascar -> Birmania
$ws.Cells.Item(138,1).Value = "Birmania"
$ws.Cells.Item(138,3).Value = 2
$ws.Cells.Item(138,4).Value = 7
$ws.Cells.Item(138,5).Value = 109
$ws.Cells.Item(138,6).Value = 0
$ws.Cells.Item(138,8).Value = 5
# Row 139: Trinidad yTobago -> Madagascar
$ws.Cells.Item(139,1).Value = "Madagascar"
$ws.Cells.Item(139,2).Value = 121
$ws.Cells.Item(139,3).Value = 0
$ws.Cells.Item(139,4).Value = 44
$ws.Cells.Item(139,5).Value = 77
$ws.Cells.Item(139,6).Value = 1
$ws.Cells.Item(139,8).Value = 0
# Row 140: Etiopia -> Trinidad yTobago
$ws.Cells.Item(140,1).Value = "Trinidad yTobago"
$ws.Cells.Item(140,2).Value = 115
$ws.Cells.Item(140,3).Value = 1
$ws.Cells.Item(140,4).Value = 28
$ws.Cells.Item(140,5).Value = 79
$ws.Cells.Item(140,8).Value = 8
# Row 141: Sudan -> Etiopia
$ws.Cells.Item(141,1).Value = "Etiopia"
$ws.Cells.Item(141,2).Value = 114
$ws.Cells.Item(141,3).Value = 3
$ws.Cells.Item(141,4).Value = 16
$ws.Cells.Item(141,5).Value = 95
$ws.Cells.Item(141,8).Value = 3
# Row 154: San Martin (Parte Holandesa) -> San Martin (Parte Holandesa)
$ws.Cells.Item(154,2).Value = 68
$ws.Cells.Item(154,3).Value = 1
$ws.Cells.Item(154,5).Value = 46
# Row 159: Haiti -> Libia
$ws.Cells.Item(159,1).Value = "Libia"
$ws.Cells.Item(159,2).Value = 59
$ws.Cells.Item(159,3).Value = 8
$ws.Cells.Item(159,4).Value = 15
$ws.Cells.Item(159,5).Value = 43
$ws.Cells.Item(159,8).Value = 1
# Row 160: Polinesia Francesa -> Haiti
$ws.Cells.Item(160,1).Value = "Haiti"
$ws.Cells.Item(160,3).Value = 0
$ws.Cells.Item(160,4).Value = 0
$ws.Cells.Item(160,5).Value = 54
$ws.Cells.Item(160,6).Value = 0
$ws.Cells.Item(160,8).Value = 3
# Row 161: Benin -> Polinesia Francesa
$ws.Cells.Item(161,1).Value = "Polinesia Francesa"
$ws.Cells.Item(161,2).Value = 57
$ws.Cells.Item(161,3).Value = 1
$ws.Cells.Item(161,4).Value = 35
$ws.Cells.Item(161,5).Value = 22
$ws.Cells.Item(161,6).Value = 1
$ws.Cells.Item(161,8).Value = 0
# Row 162: Libia -> Benin
$ws.Cells.Item(162,1).Value = "Benin"
$ws.Cells.Item(162,2).Value = 54
$ws.Cells.Item(162,4).Value = 27
$ws.Cells.Item(162,5).Value = 26
# Row 171: San Martin (Parte Francesa) -> San Martin (Parte Francesa)
$ws.Cells.Item(171,2).Value = 38
$ws.Cells.Item(171,3).Value = 1
$ws.Cells.Item(171,5).Value = 17
$ws.Cells.Item(171,6).Value = 3
# Row 182: Fiyi -> Belice
$ws.Cells.Item(182,1).Value = "Belice"
$ws.Cells.Item(182,4).Value = 2
$ws.Cells.Item(182,5).Value = 14
$ws.Cells.Item(182,6).Value = 1
$ws.Cells.Item(182,8).Value = 2
# Row 183: Belice -> Malaui
$ws.Cells.Item(183,1).Value = "Malaui"
$ws.Cells.Item(183,3).Value = 1
$ws.Cells.Item(183,4).Value = 3
$ws.Cells.Item(183,5).Value = 13
# Row 184: Malaui -> Fiyi
$ws.Cells.Item(184,1).Value = "Fiyi"
$ws.Cells.Item(184,3).Value = 0
$ws.Cells.Item(184,4).Value = 8
$ws.Cells.Item(184,5).Value = 10
$ws.Cells.Item(184,6).Value = 0
$ws.Cells.Item(184,8).Value = 0
# Row 215: San Pedro y Miquelon -> Yemen
$ws.Cells.Item(215,1).Value = "Yemen"
# Row 216: Yemen -> San Pedro y Miquelon
$ws.Cells.Item(216,1).Value = "San Pedro y Miquelon"
